$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MapSet")
$ws.Activate()

# New NPC quest rows to append below the existing table (rows 11-17).
$rows = @(
    @{ Row = 11; Id = 42030005; Name = "奥莱伊李";   Key = "npcaolai" },
    @{ Row = 12; Id = 42030006; Name = "科迪";       Key = "npckedi" },
    @{ Row = 13; Id = 42030007; Name = "威阿伊丁";   Key = "npcweia" },
    @{ Row = 14; Id = 42030008; Name = "米兰达";     Key = "npcmilanda" },
    @{ Row = 15; Id = 42030009; Name = "贝露凯伊鲁"; Key = "npcbeilukai" },
    @{ Row = 16; Id = 42030010; Name = "雷洛比克";   Key = "npcleiluo" },
    @{ Row = 17; Id = 42030011; Name = "巴鲁迪亚斯"; Key = "npcbaludi" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Id          # A - Id
    $ws.Cells.Item($rowNum, 2).Value = $r.Name         # B - Name
    $ws.Cells.Item($rowNum, 3).Value = 2                # C - Type
    $ws.Cells.Item($rowNum, 4).Value = 0                # D - Level
    $ws.Cells.Item($rowNum, 6).Value = $r.Key           # F - Ename
    $ws.Cells.Item($rowNum, 7).Value = $r.Key           # G - Figue
    $ws.Cells.Item($rowNum, 8).Value = $r.Key           # H - Script
    $ws.Cells.Item($rowNum, 9).Value = "true"           # I - TriggerMulti
}

$ws.Range("F15").Select()
